$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 23650
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 23650
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 23650
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -24586

$ws.Range("H23").Value = 23650
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 23650
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 23650
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -24118

$ws.Range("H34").Value = 883.75
$ws.Range("I34").Value = 883.75
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 883.75
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -680.75
$ws.Range("N34").ClearContents()

$ws.Range("H36").Value = 883.75
$ws.Range("I36").Value = 883.75
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 883.75
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -168.75
$ws.Range("N36").ClearContents()

$ws.Range("H70").Value = 1395.375
$ws.Range("I70").Value = 1078.5
$ws.Range("J70").Value = 1585.5
$ws.Range("K70").Value = 3235.5
$ws.Range("L70").Value = 4756.5
$ws.Range("M70").Value = -2965.5
$ws.Range("N70").Value = -5296.5

$ws.Range("H73").Value = 1395.375
$ws.Range("I73").Value = 1078.5
$ws.Range("J73").Value = 1585.5
$ws.Range("K73").Value = 3235.5
$ws.Range("L73").Value = 4756.5
$ws.Range("M73").Value = -2299.5
$ws.Range("N73").Value = -6628.5

$ws.Range("H86").Value = 33336498
$ws.Range("I86").Value = 71430950
$ws.Range("J86").Value = 3850
$ws.Range("K86").Value = 71430950
$ws.Range("L86").Value = 3850
$ws.Range("M86").Value = -71429827
$ws.Range("N86").Value = -6096

$ws.Range("H89").Value = 33336498
$ws.Range("I89").Value = 71430950
$ws.Range("J89").Value = 3850
$ws.Range("K89").Value = 357154750
$ws.Range("L89").Value = 19250
$ws.Range("M89").Value = -357149134
$ws.Range("N89").Value = -30482

$ws.Range("H132").Value = 3409.5144
$ws.Range("I132").Value = 3358.9688
$ws.Range("J132").Value = 3948.6667
$ws.Range("K132").Value = 10076.9064
$ws.Range("L132").Value = 11846.0001
$ws.Range("M132").Value = -7546.9064
$ws.Range("N132").Value = -16906.0001

$ws.Range("H135").Value = 4449.6
$ws.Range("I135").Value = 3828
$ws.Range("J135").Value = 5900
$ws.Range("K135").Value = 34452
$ws.Range("L135").Value = 53100
$ws.Range("M135").Value = -31917
$ws.Range("N135").Value = -58170

$ws.Range("H137").Value = 1185.8182
$ws.Range("I137").Value = 890.7037
$ws.Range("J137").Value = 2513.8333
$ws.Range("K137").Value = 2672.1111
$ws.Range("L137").Value = 7541.499899999999
$ws.Range("M137").Value = -122.1111000000001
$ws.Range("N137").Value = -12641.4999

$ws.Range("H138").Value = 2358.125
$ws.Range("I138").Value = 1847.5
$ws.Range("J138").Value = 2868.75
$ws.Range("K138").Value = 5542.5
$ws.Range("L138").Value = 8606.25
$ws.Range("M138").Value = -402.5
$ws.Range("N138").Value = -18886.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H39").Value = 15890.6
$ws.Range("J39").Value = 15890.6
$ws.Range("L39").Value = 15890.6
$ws.Range("N39").Value = -16668.6

$ws.Range("H112").Value = 34469
$ws.Range("J112").Value = 34469
$ws.Range("L112").Value = 34469
$ws.Range("N112").Value = -37423

$ws.Range("H134").Value = 2626.585
$ws.Range("I134").Value = 1477.5294
$ws.Range("J134").Value = 4682.7896
$ws.Range("K134").Value = 4432.5882
$ws.Range("L134").Value = 14048.3688
$ws.Range("M134").Value = -1897.5882
$ws.Range("N134").Value = -19118.3688

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 22794.715
$ws.Range("I26").Value = 2375
$ws.Range("J26").Value = 50021
$ws.Range("K26").Value = 2375
$ws.Range("L26").Value = 50021
$ws.Range("M26").Value = -2088
$ws.Range("N26").Value = -50595

$ws.Range("H31").Value = 6767.522
$ws.Range("I31").Value = 8509
$ws.Range("J31").Value = 1833.3334
$ws.Range("K31").Value = 8509
$ws.Range("L31").Value = 1833.3334
$ws.Range("M31").Value = -8214
$ws.Range("N31").Value = -2423.3334

$ws.Range("H34").Value = 6767.522
$ws.Range("I34").Value = 8509
$ws.Range("J34").Value = 1833.3334
$ws.Range("K34").Value = 8509
$ws.Range("L34").Value = 1833.3334
$ws.Range("M34").Value = -8307
$ws.Range("N34").Value = -2237.3334

$ws.Range("H44").Value = 96266
$ws.Range("I44").Value = 118354.664
$ws.Range("K44").Value = 118354.664
$ws.Range("M44").Value = -117912.664

$ws.Range("H54").Value = 29962.666
$ws.Range("J54").Value = 29962.666
$ws.Range("L54").Value = 29962.666
$ws.Range("N54").Value = -31278.666

$ws.Range("H132").Value = 10037.571
$ws.Range("I132").Value = 11551.5
$ws.Range("K132").Value = 34654.5
$ws.Range("M132").Value = -32124.5

$ws.Range("H134").Value = 2000.7273
$ws.Range("I134").Value = 1926.8
$ws.Range("J134").Value = 2740
$ws.Range("K134").Value = 5780.4
$ws.Range("L134").Value = 8220
$ws.Range("M134").Value = -3245.4
$ws.Range("N134").Value = -13290

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 611.8333
$ws.Range("I107").Value = 257.5
$ws.Range("J107").Value = 682.7
$ws.Range("K107").Value = 772.5
$ws.Range("L107").Value = 2048.1
$ws.Range("M107").Value = 1147.5
$ws.Range("N107").Value = -5888.1

$ws.Range("H113").Value = 1586.2222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5459.2
$ws.Range("I40").Value = 5145.6
$ws.Range("J40").Value = 6400
$ws.Range("K40").Value = 5145.6
$ws.Range("L40").Value = 6400
$ws.Range("M40").Value = -5009.6
$ws.Range("N40").Value = -6672

$ws.Range("H64").Value = 22575
$ws.Range("J64").Value = 22575
$ws.Range("L64").Value = 22575
$ws.Range("N64").Value = -23025

$ws.Range("H67").Value = 22575
$ws.Range("J67").Value = 22575
$ws.Range("L67").Value = 22575
$ws.Range("N67").Value = -24135

$ws.Range("H110").Value = 28633
$ws.Range("J110").Value = 28633
$ws.Range("L110").Value = 28633
$ws.Range("N110").Value = -36813

$ws.Range("H132").Value = 13896951
$ws.Range("I132").Value = 4296.263
$ws.Range("J132").Value = 29424036
$ws.Range("K132").Value = 12888.789
$ws.Range("L132").Value = 88272108
$ws.Range("M132").Value = -10358.789
$ws.Range("N132").Value = -88277168

$ws.Range("H135").Value = 104987
$ws.Range("J135").Value = 104987
$ws.Range("L135").Value = 104987
$ws.Range("N135").Value = -115127

$ws.Range("H136").Value = 11387.692
$ws.Range("I136").Value = 2503.889
$ws.Range("K136").Value = 7511.667
$ws.Range("M136").Value = -4961.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 38766.668
$ws.Range("J42").Value = 38650
$ws.Range("L42").Value = 38650
$ws.Range("N42").Value = -39406

$ws.Range("H136").Value = 1421.2084
$ws.Range("I136").Value = 1437.75
$ws.Range("J136").Value = 1388.125
$ws.Range("K136").Value = 4313.25
$ws.Range("L136").Value = 4164.375
$ws.Range("M136").Value = -1763.25
$ws.Range("N136").Value = -9264.375
